$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Varsha", "Singhal", "MindTree"),
    @("Nidhi", "Choudhary", "Infosys"),
    @("Upasana", "Sinha", "Cognizant"),
    @("Ruchita", "Kadam", "IBM")
)

$row = 5
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

$ws.Range("C8").Select()
